# Add a new "Server 7" row to the servers table.
#
# Mirrors what Excel recorded when the author appended a row below the
# existing data (Server6 / 192.168.7.1 / Server3 block ends at row 10):
#   A11 = "Server 7"
#   B11 = "192.168.9.1"
#   C11 = "Server6"          (ConnectsTo -> existing "Server6" entry)
#
# The new A11:B11 cells pick up the same formatting (vertical-center +
# wrap text) already used by the rest of the data rows, so the format is
# copied from the row directly above (A10:B10) rather than re-describing
# it by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the look of the row above (vertical-center + wrap alignment) onto
# the new row's Name/IP cells before filling in values.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Server 7"
$ws.Range("B11").Value = "192.168.9.1"
$ws.Range("C11").Value = "Server6"

# Match Excel's recorded post-edit selection (cell C11).
$ws.Range("C11").Select() | Out-Null
